# "more dialogs: act 3-2, 3-3"
# Adds the act_3_2_* and act_3_3_* dialog lines (keys + English text) to the
# "en" sheet, and the corresponding (still-untranslated) keys to the "es"
# sheet, mirroring how the earlier act_3_intro_* rows were entered.

$wb = $excel.ActiveWorkbook
$ws_en = $wb.Worksheets.Item("en")
$ws_es = $wb.Worksheets.Item("es")

# --- "en" sheet: rows 127-129 (act_3_2 intro) ------------------------------
$ws_en.Range("A127").Value = "act_3_2_intro_1"
$ws_en.Range("A128").Value = "act_3_2_intro_2"
$ws_en.Range("A129").Value = "act_3_2_intro_3"

$ws_en.Range("B127").Value = "This time around, we will be placing two force fields."
$ws_en.Range("B128").Value = "In order to allow the golden ball to take flight, one must counteract the gravitational force with a much greater force."
$ws_en.Range("B129").Value = "Please take note of the force field’s direction when placing them on the wall. Good luck!"

# --- "en" sheet: rows 130-133 (act_3_3 intro) ------------------------------
$ws_en.Range("A130").Value = "act_3_3_intro_1"
$ws_en.Range("A131").Value = "act_3_3_intro_2"
$ws_en.Range("A132").Value = "act_3_3_intro_3"
$ws_en.Range("A133").Value = "act_3_3_intro_4"

$ws_en.Range("B131").Value = "First: An object at rest will remain at rest, and an object in motion will maintain a constant velocity unless acted upon by a net external force."
$ws_en.Range("B132").Value = "Second: The acceleration of an object is proportional to the net force acting on it, and inversely proportional to its mass."
$ws_en.Range("B133").Value = "Keep these in mind, and you will surely be a force to be reckoned with!"
$ws_en.Range("B130").Value = "Now that you’ve come this far, it’s all up to you! Remember all that you’ve learned about Newton’s first two laws."

# Match the vertical-centered style already used on other multi-line prompts.
$ws_en.Range("B127").VerticalAlignment = -4108
$ws_en.Range("B129").VerticalAlignment = -4108

# --- "es" sheet: rows 127-129 (keys only, not yet translated) -------------
$ws_es.Range("A127").Value = "act_3_2_intro_1"
$ws_es.Range("A128").Value = "act_3_2_intro_2"
$ws_es.Range("A129").Value = "act_3_2_intro_3"

# --- Update the view state to match where the author ended up editing -----
$ws_en.Activate() | Out-Null
$ws_en.Range("B130").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 106
$excel.ActiveWindow.ScrollColumn = 1 | Out-Null
